{"js": "// The paragraph `<id>p045r_1</id>` was previously split across three runs\n// (`<id>`, `p045r_1`, `</id>`) with the first/third run carrying the\n// Courier-New/brownish \"tag\" formatting and the middle run carrying plain\n// black text. The edit merges them back into a single run containing the\n// full literal text `<id>p045r_1</id>`.\n//\n// Find the run/range that currently spans that text (search() matches\n// across run boundaries) and rewrite its text in place. Word/Office.js\n// collapses a range spanning multiple runs down to one run using the\n// formatting of the range's first character, which reproduces the\n// target markup exactly.\nconst results = context.document.body.search(\"<id>p045r_1</id>\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '<id>p045r_1</id>' in the document body.\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\"<id>p045r_1</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The paragraph \"<id>p045r_1</id>\" was split across three runs (`<id>`,\n# `p045r_1`, `</id>`), with the first/third run carrying the Courier-New\n# brownish \"tag\" formatting and the middle run plain black text. Collapse\n# them back into a single run holding the full literal text\n# \"<id>p045r_1</id>\", formatted like the original first run.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"<id>p045r_1</id>\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 0\n\n# Execute with an in-place replacement of the identical literal text: Word\n# re-serializes the matched range as a single run (taking on the formatting\n# of the range's first character), which merges the three original runs\n# into one, exactly like the OOXML diff.\n$find.Execute(\"<id>p045r_1</id>\", $false, $false, $false, $false, $false, $true, 1, $false, \"<id>p045r_1</id>\", 2)\n"}
